$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # 1. Update "想去人数" (F column) counts that increased.
    $ws.Range("F2").Value = 46
    $ws.Range("F3").Value = 3070
    $ws.Range("F5").Value = 160
    $ws.Range("F7").Value = 1692
    $ws.Range("F9").Value = 89
    $ws.Range("F12").Value = 1387
    $ws.Range("F14").Value = 530
    $ws.Range("F16").Value = 43
    $ws.Range("F20").Value = 125
    $ws.Range("F22").Value = 107
    $ws.Range("F23").Value = 3243
    $ws.Range("F24").Value = 394
    $ws.Range("F25").Value = 149
    $ws.Range("F26").Value = 335

    # 2. Insert a new row 29 for the newly-announced event, pushing the old
    #    row 29 (南昌·代号鸢盛花行only) down to row 30.
    $ws.Rows.Item(29).Insert()

    # Re-create row 29's formatting (style) by copying the row directly
    # above it, then overwrite with the new event's data.
    $ws.Range("A28").Copy($ws.Range("A29"))

    $ws.Range("A29").Value = 28

    # Force the date-looking text into a real text cell (otherwise Excel's
    # smart typing would reinterpret it as a date serial number), then
    # drop the leftover "text" number format again so the cell matches its
    # plain, unformatted siblings.
    $ws.Range("B29").NumberFormat = "@"
    $ws.Range("B29").Value = "2024-05-02"
    $ws.Range("B29").ClearFormats()

    $ws.Range("C29").Value = "江西·ShiningStaR数字互娱嘉年华"
    $ws.Range("D29").Value = "前湖大道欣悦湖体育馆 欣悦湖体育馆"
    $ws.Range("E29").Value = "2024.05.02 09:30-05.04 17:00"
    $ws.Range("F29").Value = 73
    $ws.Range("G29").Value = 55
    $ws.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=83180"
    $ws.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202403/EqmGU5NC1711015780862.jpeg"

    # The event that used to be row 29 is now row 30; its sequence number
    # (column A) increments from 28 to 29.
    $ws.Range("A30").Value = 29
}
